# Applies the "Updated symbol list" edit: refresh crypto price/volume-label cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.74"
$ws.Range("D2").ClearFormats()

$ws.Range("D4").Value = "'5.411"
$ws.Range("D4").ClearFormats()

$ws.Range("D5").Value = "'0.05799"
$ws.Range("D5").ClearFormats()

$ws.Range("D7").Value = "'6.330"
$ws.Range("D7").ClearFormats()

$ws.Range("D8").Value = "'0.8069"
$ws.Range("D8").ClearFormats()

$ws.Range("D9").Value = "'0.9620"
$ws.Range("D9").ClearFormats()

$ws.Range("E9").Value = "8FTXTokenFTTBestin24h"

$ws.Range("D10").Value = "'0.1427"
$ws.Range("D10").ClearFormats()

$ws.Range("D11").Value = "'0.07516"
$ws.Range("D11").ClearFormats()

$ws.Range("D12").Value = "'0.03222"
$ws.Range("D12").ClearFormats()

$ws.Range("D13").Value = "'0.03022"
$ws.Range("D13").ClearFormats()

$ws.Range("D14").Value = "'4.146"
$ws.Range("D14").ClearFormats()

$ws.Range("D15").Value = "'0.09405"
$ws.Range("D15").ClearFormats()

$ws.Range("D16").Value = "'0.001586"
$ws.Range("D16").ClearFormats()

$ws.Range("D17").Value = "'0.04807"
$ws.Range("D17").ClearFormats()

$ws.Range("D18").Value = "'0.0005898"
$ws.Range("D18").ClearFormats()

$ws.Range("D19").Value = "'0.005313"
$ws.Range("D19").ClearFormats()

$ws.Range("D20").Value = "'0.004110"
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = "19HotbitTokenHTBWorstin24h"

$ws.Range("D21").Value = "'0.0009978"
$ws.Range("D21").ClearFormats()

$ws.Range("D23").Value = "'3.739"
$ws.Range("D23").ClearFormats()

$ws.Range("D24").Value = "'2.237"
$ws.Range("D24").ClearFormats()

$ws.Range("D26").Value = "'0.1258"
$ws.Range("D26").ClearFormats()

$ws.Range("D27").Value = "'0.0003123"
$ws.Range("D27").ClearFormats()

$ws.Range("D40").Value = "'0.03895"
$ws.Range("D40").ClearFormats()

$ws.Range("D41").Value = "'0.006359"
$ws.Range("D41").ClearFormats()

$ws.Range("D42").Value = "'0.1076"
$ws.Range("D42").ClearFormats()

$ws.Range("D43").Value = "'0.002660"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.006705"
$ws.Range("D44").ClearFormats()

$ws.Range("D45").Value = "'0.00005590"
$ws.Range("D45").ClearFormats()

$ws.Range("D47").Value = "'0.3899"
$ws.Range("D47").ClearFormats()

$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

$ws.Range("D48").Value = "'0.1470"
$ws.Range("D48").ClearFormats()

$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").ClearFormats()
